$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.097.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.06%  "
$ws.Range("E3").Value = "  -3.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4001"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08380"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.042"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.882.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.39%  "
$ws.Range("E14").Value = "  -5.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.044"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001064"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06587"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.737"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.090.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.306"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.140.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("E28").Value = "  -2.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.740"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.127"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9734"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.14%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09630"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.459"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.631"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.543"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.279"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.778"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02290"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06142"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6145"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1905"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("E45").Value = "  -3.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5858"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("E48").Value = "  -5.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.437"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06911"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.24%  "
